$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) <sup> -> <corr>  and  </sup> -> </corr>
# ------------------------------------------------------------------
$d.Content.Find.Execute("<sup>", $true, $false, $false, $false, $false, $true, 1, $false, "<corr>", 2) | Out-Null
$d.Content.Find.Execute("</sup>", $true, $false, $false, $false, $false, $true, 1, $false, "</corr>", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "A bit of " -> "A " + "little " (as two separate runs; the new
#    "little " run must carry the same "no explicit color" formatting
#    as the existing "left-middle" run, not the black color baked
#    into "A bit of "). We get that formatting by copy/pasting the
#    "left-middle" run (non-destructive - source stays untouched)
#    and then overwriting the pasted text.
# ------------------------------------------------------------------
$fmtSrc = $d.Content.Duplicate
$fmtSrc.Find.Execute("left-middle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fmtSrc.Copy()

$target = $d.Content.Duplicate
$target.Find.Execute("A bit of ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Text = "A "
$target.Collapse(0)
$insertStart = $target.Start
$target.Paste()

$pasted = $d.Range($insertStart, $insertStart + 11)
$pasted.Text = "little "

# ------------------------------------------------------------------
# 3) remove the "," run right after the first </m> ("tallow</m>,")
# ------------------------------------------------------------------
$commaFind = $d.Content.Duplicate
$commaFind.Find.Execute("tallow</m>,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$comma = $d.Range($commaFind.End - 1, $commaFind.End)
$comma.Text = ""

Write-Output "done"
